# Sandbox 7: Player and weapon classes under construction.
# Game now resets and has a score; recreates the managers and player.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Balloon Sandbox")

# ------------------------------------------------------------------
# 1) Rename the sandbox tab: "Balloon Sandbox" -> "Game Sandbox"
# ------------------------------------------------------------------
$ws.Name = "Game Sandbox"

# ------------------------------------------------------------------
# 1b) Extend rows 7-12 with a new column T (all "Bad" style, matching
#     the rest of those rows' S column).
# ------------------------------------------------------------------
$badSrcForT = $ws.Cells.Item(4, 2)  # B4 - Bad style, empty
for ($row = 7; $row -le 12; $row++) {
    $badSrcForT.Copy($ws.Cells.Item($row, 20))
}

# ------------------------------------------------------------------
# 2) Add row 14 (task #8) BEFORE touching row 13, reusing row 13's
#    current (untouched, all "Normal w/ border" style) cells as the
#    format source so no new cell-style entries are created.
# ------------------------------------------------------------------
for ($col = 2; $col -le 19; $col++) {
    $src = $ws.Cells.Item(13, $col)
    $dst = $ws.Cells.Item(14, $col)
    $src.Copy($dst)
}
# Column T (20) is brand new - any still-untouched style-13 cell works as source.
$ws.Cells.Item(13, 2).Copy($ws.Cells.Item(14, 20))

# A14: task number 8, same centered style as A13.
$ws.Cells.Item(13, 1).Copy($ws.Cells.Item(14, 1))
$ws.Cells.Item(14, 1).Value = 8

# ------------------------------------------------------------------
# 3) Row 13 (task #7) moves from "empty/unset" look to the same
#    filled-in Neutral/Bad pattern used by the other data rows.
# ------------------------------------------------------------------
$neutralSrc = $ws.Cells.Item(3, 2)  # B3 - Neutral style, empty
$badSrc     = $ws.Cells.Item(4, 2)  # B4 - Bad style, empty

$row13Styles = @{
    2  = "Neutral"  # B13
    3  = "Neutral"  # C13
    4  = "Neutral"  # D13
    5  = "Bad"      # E13
    6  = "Bad"      # F13
    7  = "Bad"      # G13
    8  = "Neutral"  # H13
    9  = "Bad"      # I13
    10 = "Neutral"  # J13
    11 = "Bad"      # K13
    12 = "Bad"      # L13
    13 = "Neutral"  # M13
    14 = "Neutral"  # N13
    15 = "Neutral"  # O13
    16 = "Neutral"  # P13
    17 = "Bad"      # Q13
    18 = "Neutral"  # R13
    19 = "Neutral"  # S13
    20 = "Neutral"  # T13
}

foreach ($col in $row13Styles.Keys) {
    $dst = $ws.Cells.Item(13, $col)
    if ($row13Styles[$col] -eq "Neutral") {
        $neutralSrc.Copy($dst)
    } else {
        $badSrc.Copy($dst)
    }
}

# ------------------------------------------------------------------
# 4) Header row 6: rework the "Shotgun"/"Balloon Manager" columns
#    into the new "Character Managers" / "Player" / "Weapon" trio.
# ------------------------------------------------------------------
$r6 = $ws.Cells.Item(6, 18)  # was "Balloon Manager"
$s6 = $ws.Cells.Item(6, 19)  # was "Shotgun"
$t6 = $ws.Cells.Item(6, 20)  # brand new column

# Give the new T6 header the same header formatting as its neighbour
# before overwriting R6/S6 text (R6 still carries the original style).
$r6.Copy($t6)

$r6.Value = "Character Managers"
$s6.Value = "Player"
$t6.Value = "Weapon"

# ------------------------------------------------------------------
# 5) Selection moves from R6 to P10.
# ------------------------------------------------------------------
$ws.Activate()
$ws.Range("P10").Select()
